$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 with new (re-matched) evaluation values
$ws.Range("B2").Value = 0.2247871045231402
$ws.Range("C2").Value = 2.035146975718185
$ws.Range("D2").Value = 17.83605474045861
$ws.Range("E2").Value = 4.223275356930757
$ws.Range("F2").Value = 4.269679570594401
$ws.Range("G2").Value = 41

$ws.Range("B3").Value = 0.2728254128479093
$ws.Range("C3").Value = 1.970770711136308
$ws.Range("D3").Value = 15.69395479060738
$ws.Range("E3").Value = 3.961559641177623
$ws.Range("F3").Value = 3.967096014706342
$ws.Range("G3").Value = 133

# Add new row 4 ("Q1") - copy formatting from row 3's label cell first
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = "Q1"

$ws.Range("B4").Value = 0.171141602336921
$ws.Range("C4").Value = 1.274244979641737
$ws.Range("D4").Value = 5.743778122149821
$ws.Range("E4").Value = 2.396618059297272
$ws.Range("F4").Value = 2.408541433280348
$ws.Range("G4").Value = 67
